# Update "想去人数" (interested-people count) figures in column F across sheets,
# reflecting newly generated output data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3113
$ws1.Range("F5").Value = 6833
$ws1.Range("F6").Value = 1888
$ws1.Range("F8").Value = 59
$ws1.Range("F12").Value = 12
$ws1.Range("F14").Value = 160
$ws1.Range("F15").Value = 30

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 9

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3113
$ws4.Range("F3").Value = 9
$ws4.Range("F6").Value = 6833
$ws4.Range("F7").Value = 1888
$ws4.Range("F9").Value = 59
$ws4.Range("F13").Value = 12
$ws4.Range("F15").Value = 160
$ws4.Range("F16").Value = 30
